# Auto-generated: apply 2024-12-29 CTA violent-crime-ytd data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("F2").Value = 98  # 97 -> 98
$ws.Range("J2").Value = 130  # 129 -> 130
$ws.Range("B3").Value = 82  # 81 -> 82
$ws.Range("C3").Value = 82  # 81 -> 82
$ws.Range("F3").Value = 146  # 145 -> 146
$ws.Range("I3").Value = 199  # 198 -> 199
$ws.Range("K3").Value = 232  # 230 -> 232
$ws.Range("B6").Value = 395  # 391 -> 395
$ws.Range("C6").Value = 508  # 507 -> 508
$ws.Range("D6").Value = 443  # 440 -> 443
$ws.Range("E6").Value = 512  # 509 -> 512
$ws.Range("F6").Value = 579  # 575 -> 579
$ws.Range("H6").Value = 480  # 479 -> 480
$ws.Range("J6").Value = 437  # 435 -> 437
$ws.Range("B7").Value = 535  # 530 -> 535
$ws.Range("C7").Value = 670  # 668 -> 670
$ws.Range("D7").Value = 686  # 683 -> 686
$ws.Range("E7").Value = 750  # 747 -> 750
$ws.Range("F7").Value = 833  # 827 -> 833
$ws.Range("H7").Value = 769  # 768 -> 769
$ws.Range("I7").Value = 866  # 865 -> 866
$ws.Range("J7").Value = 831  # 828 -> 831
$ws.Range("K7").Value = 948  # 946 -> 948

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("F3").Value = 5  # 4 -> 5
$ws.Range("E6").Value = 56  # 55 -> 56
$ws.Range("E7").Value = 69  # 68 -> 69
$ws.Range("F7").Value = 64  # 63 -> 64

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("F5").Value = 14  # 12 -> 14
$ws.Range("F6").Value = 18  # 16 -> 18

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("B4").Value = 13  # 12 -> 13
$ws.Range("D4").Value = 5  # 4 -> 5
$ws.Range("B5").Value = 17  # 16 -> 17
$ws.Range("D5").Value = 6  # 5 -> 6

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J5").Value = 5  # 4 -> 5
$ws.Range("J6").Value = 16  # 15 -> 16

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("B6").Value = 34  # 33 -> 34
$ws.Range("B7").Value = 40  # 39 -> 40

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("B4").Value = 11  # 9 -> 11
$ws.Range("B5").Value = 15  # 13 -> 15

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F5").Value = 18  # 16 -> 18
$ws.Range("C19").Value = 11  # 10 -> 11
$ws.Range("F19").Value = 26  # 25 -> 26
$ws.Range("E20").Value = 5  # 4 -> 5
$ws.Range("B28").Value = 40  # 39 -> 40
$ws.Range("E32").Value = 69  # 68 -> 69
$ws.Range("F32").Value = 64  # 63 -> 64
$ws.Range("B45").Value = 4  # 3 -> 4
$ws.Range("K47").Value = 24  # 23 -> 24
$ws.Range("C53").Value = 62  # 61 -> 62
$ws.Range("D53").Value = 79  # 78 -> 79
$ws.Range("F53").Value = 88  # 87 -> 88
$ws.Range("I53").Value = 129  # 128 -> 129
$ws.Range("J53").Value = 130  # 128 -> 130
$ws.Range("H61").Value = 9  # 8 -> 9
$ws.Range("K63").Value = 10  # 9 -> 10
$ws.Range("E65").Value = 20  # 19 -> 20
$ws.Range("F74").Value = 12  # 11 -> 12
$ws.Range("D79").Value = 7  # 6 -> 7
$ws.Range("B80").Value = 17  # 16 -> 17
$ws.Range("D80").Value = 6  # 5 -> 6
$ws.Range("B82").Value = 15  # 13 -> 15
$ws.Range("J86").Value = 16  # 15 -> 16
$ws.Range("B98").Value = 535  # 530 -> 535
$ws.Range("C98").Value = 670  # 668 -> 670
$ws.Range("D98").Value = 686  # 683 -> 686
$ws.Range("E98").Value = 750  # 747 -> 750
$ws.Range("F98").Value = 833  # 827 -> 833
$ws.Range("H98").Value = 769  # 768 -> 769
$ws.Range("I98").Value = 866  # 865 -> 866
$ws.Range("J98").Value = 831  # 828 -> 831
$ws.Range("K98").Value = 948  # 946 -> 948

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("F2").Value = 8  # 7 -> 8
$ws.Range("J2").Value = 22  # 21 -> 22
$ws.Range("I3").Value = 31  # 30 -> 31
$ws.Range("C6").Value = 44  # 43 -> 44
$ws.Range("D6").Value = 49  # 48 -> 49
$ws.Range("J6").Value = 63  # 62 -> 63
$ws.Range("C7").Value = 62  # 61 -> 62
$ws.Range("D7").Value = 79  # 78 -> 79
$ws.Range("F7").Value = 88  # 87 -> 88
$ws.Range("I7").Value = 129  # 128 -> 129
$ws.Range("J7").Value = 130  # 128 -> 130

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("D5").Value = 5  # 4 -> 5
$ws.Range("D6").Value = 7  # 6 -> 7

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("E5").Value = 17  # 16 -> 17
$ws.Range("E6").Value = 20  # 19 -> 20

$ws = $wb.Worksheets.Item("River North")
$ws.Range("F5").Value = 11  # 10 -> 11
$ws.Range("F6").Value = 12  # 11 -> 12

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 5  # 4 -> 5
$ws.Range("J5").Value = 10  # 9 -> 10

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("C3").Value = 1  # new cell = 1
$ws.Range("F6").Value = 19  # 18 -> 19
$ws.Range("C7").Value = 11  # 10 -> 11
$ws.Range("F7").Value = 26  # 25 -> 26

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("B3").Value = 1  # new cell = 1

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("B6").Value = 4  # 3 -> 4

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 5  # 4 -> 5
$ws.Range("K6").Value = 24  # 23 -> 24

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("E5").Value = 4  # 3 -> 4
$ws.Range("E6").Value = 5  # 4 -> 5
